$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$updates = @{
    5  = -5
    6  = 1
    15 = -1
    20 = -5
    25 = -1
    27 = -1
    28 = 0
    29 = 3
    31 = -1
    37 = -1
    41 = -1
    46 = 2
    50 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
